$wb = $excel.ActiveWorkbook

# Grab an existing header cell to copy its formatting (bold, centered, bordered)
$ws1 = $wb.Worksheets.Item(1)
$srcHeaderRange = $ws1.Range("A1:D1")

# Add the new worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = "'" + $headers[$i]
}

# Copy header formatting (bold font, centered alignment, thin border) onto the new header row
$srcHeaderRange.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows
$data = @(
    @("4402", "", "", "", "", "NO"),
    @("4406", "", "", "", "", "NO"),
    @("4745", "", "", "", "", "NO")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($val -ne "") {
            $newSheet.Cells.Item($r + 2, $c + 1).Value = "'" + $val
        }
    }
}

$newSheet.Range("A1").Select() | Out-Null
